$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value2 = 30777
$ws.Range("I7").Value2 = 0
$ws.Range("J7").Value2 = 30777
$ws.Range("K7").Value2 = 0
$ws.Range("L7").Value2 = 30777
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value2 = -31001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value2 = 30777
$ws.Range("I14").Value2 = 0
$ws.Range("J14").Value2 = 30777
$ws.Range("K14").Value2 = 0
$ws.Range("L14").Value2 = 30777
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value2 = -31159

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value2 = 432.33334
$ws.Range("I19").Value2 = 500
$ws.Range("J19").Value2 = 398.5
$ws.Range("K19").Value2 = 500
$ws.Range("L19").Value2 = 398.5
$ws.Range("M19").Value2 = -325
$ws.Range("N19").Value2 = -748.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 5212.2856
$ws.Range("I62").Value2 = 4496.875
$ws.Range("J62").Value2 = 6166.1665
$ws.Range("K62").Value2 = 4496.875
$ws.Range("L62").Value2 = 6166.1665
$ws.Range("M62").Value2 = -3872.875
$ws.Range("N62").Value2 = -7414.1665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value2 = 5212.2856
$ws.Range("I65").Value2 = 4496.875
$ws.Range("J65").Value2 = 6166.1665
$ws.Range("K65").Value2 = 22484.375
$ws.Range("L65").Value2 = 30830.8325
$ws.Range("M65").Value2 = -19364.375
$ws.Range("N65").Value2 = -37070.8325

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value2 = 773.0714
$ws.Range("I88").Value2 = 632
$ws.Range("J88").Value2 = 811.5454999999999
$ws.Range("K88").Value2 = 632
$ws.Range("L88").Value2 = 811.5454999999999
$ws.Range("M88").Value2 = -226
$ws.Range("N88").Value2 = -1623.5455

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value2 = 773.0714
$ws.Range("I91").Value2 = 632
$ws.Range("J91").Value2 = 811.5454999999999
$ws.Range("K91").Value2 = 632
$ws.Range("L91").Value2 = 811.5454999999999
$ws.Range("M91").Value2 = 772
$ws.Range("N91").Value2 = -3619.5455

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value2 = 1284.4231
$ws.Range("I137").Value2 = 1172.409
$ws.Range("J137").Value2 = 1900.5
$ws.Range("K137").Value2 = 3517.227
$ws.Range("L137").Value2 = 5701.5
$ws.Range("M137").Value2 = -967.2270000000003
$ws.Range("N137").Value2 = -10801.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value2 = 646.75
$ws.Range("I5").Value2 = 528.3333
$ws.Range("J5").Value2 = 1002
$ws.Range("K5").Value2 = 528.3333
$ws.Range("L5").Value2 = 1002
$ws.Range("M5").Value2 = -416.3333
$ws.Range("N5").Value2 = -1226

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 4225.1025
$ws.Range("I32").Value2 = 4225.1025
$ws.Range("J32").Value2 = 0
$ws.Range("K32").Value2 = 4225.1025
$ws.Range("L32").Value2 = 0
$ws.Range("M32").Value2 = -3938.1025
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value2 = 3819.3333
$ws.Range("I74").Value2 = 3921.875
$ws.Range("J74").Value2 = 2999
$ws.Range("K74").Value2 = 3921.875
$ws.Range("L74").Value2 = 2999
$ws.Range("M74").Value2 = -3047.875
$ws.Range("N74").Value2 = -4747

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value2 = 3819.3333
$ws.Range("I77").Value2 = 3921.875
$ws.Range("J77").Value2 = 2999
$ws.Range("K77").Value2 = 19609.375
$ws.Range("L77").Value2 = 14995
$ws.Range("M77").Value2 = -15241.375
$ws.Range("N77").Value2 = -23731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value2 = 740
$ws.Range("I102").Value2 = 740
$ws.Range("J102").Value2 = 0
$ws.Range("K102").Value2 = 740
$ws.Range("L102").Value2 = 0
$ws.Range("M102").Value2 = 882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value2 = 646.75
$ws.Range("I4").Value2 = 528.3333
$ws.Range("J4").Value2 = 1002
$ws.Range("K4").Value2 = 528.3333
$ws.Range("L4").Value2 = 1002
$ws.Range("M4").Value2 = -413.3333
$ws.Range("N4").Value2 = -1232

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value2 = 199.5
$ws.Range("I29").Value2 = 199.5
$ws.Range("J29").Value2 = 0
$ws.Range("K29").Value2 = 199.5
$ws.Range("L29").Value2 = 0
$ws.Range("M29").Value2 = 89.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value2 = 28767.375
$ws.Range("I82").Value2 = 10046.333
$ws.Range("J82").Value2 = 40000
$ws.Range("K82").Value2 = 10046.333
$ws.Range("L82").Value2 = 40000
$ws.Range("M82").Value2 = -9663.333000000001
$ws.Range("N82").Value2 = -40766

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value2 = 28767.375
$ws.Range("I85").Value2 = 10046.333
$ws.Range("J85").Value2 = 40000
$ws.Range("K85").Value2 = 10046.333
$ws.Range("L85").Value2 = 40000
$ws.Range("M85").Value2 = -8720.333000000001
$ws.Range("N85").Value2 = -42652

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 3961.75
$ws.Range("I86").Value2 = 3837.6
$ws.Range("J86").Value2 = 4168.6665
$ws.Range("K86").Value2 = 3837.6
$ws.Range("L86").Value2 = 4168.6665
$ws.Range("M86").Value2 = -2714.6
$ws.Range("N86").Value2 = -6414.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value2 = 3961.75
$ws.Range("I89").Value2 = 3837.6
$ws.Range("J89").Value2 = 4168.6665
$ws.Range("K89").Value2 = 19188
$ws.Range("L89").Value2 = 20843.3325
$ws.Range("M89").Value2 = -13572
$ws.Range("N89").Value2 = -32075.3325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 102.2
$ws.Range("I7").Value2 = 106.181816
$ws.Range("J7").Value2 = 91.25
$ws.Range("K7").Value2 = 106.181816
$ws.Range("L7").Value2 = 91.25
$ws.Range("M7").Value2 = 6.818184000000002
$ws.Range("N7").Value2 = -317.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1764.2142
$ws.Range("I31").Value2 = 1891
$ws.Range("J31").Value2 = 1299.3334
$ws.Range("K31").Value2 = 1891
$ws.Range("L31").Value2 = 1299.3334
$ws.Range("M31").Value2 = -1596
$ws.Range("N31").Value2 = -1889.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value2 = 1764.2142
$ws.Range("I34").Value2 = 1891
$ws.Range("J34").Value2 = 1299.3334
$ws.Range("K34").Value2 = 1891
$ws.Range("L34").Value2 = 1299.3334
$ws.Range("M34").Value2 = -1689
$ws.Range("N34").Value2 = -1703.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value2 = 38888.332
$ws.Range("I74").Value2 = 0
$ws.Range("J74").Value2 = 38888.332
$ws.Range("K74").Value2 = 0
$ws.Range("L74").Value2 = 38888.332
$ws.Range("N74").Value2 = -40636.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value2 = 38888.332
$ws.Range("I77").Value2 = 0
$ws.Range("J77").Value2 = 38888.332
$ws.Range("K77").Value2 = 0
$ws.Range("L77").Value2 = 116664.996
$ws.Range("N77").Value2 = -125400.996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value2 = 2558.8
$ws.Range("I94").Value2 = 2648.6667
$ws.Range("J94").Value2 = 2424
$ws.Range("K94").Value2 = 2648.6667
$ws.Range("L94").Value2 = 2424
$ws.Range("M94").Value2 = -2197.6667
$ws.Range("N94").Value2 = -3326

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value2 = 948.5
$ws.Range("I105").Value2 = 852.9091
$ws.Range("J105").Value2 = 2000
$ws.Range("K105").Value2 = 852.9091
$ws.Range("L105").Value2 = 2000
$ws.Range("M105").Value2 = 894.0909
$ws.Range("N105").Value2 = -5494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value2 = 0
$ws.Range("I141").Value2 = 0
$ws.Range("J141").Value2 = 0
$ws.Range("K141").Value2 = 0
$ws.Range("L141").Value2 = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value2 = 99500
$ws.Range("I37").Value2 = 0
$ws.Range("J37").Value2 = 99500
$ws.Range("K37").Value2 = 0
$ws.Range("L37").Value2 = 298500
$ws.Range("N37").Value2 = -298724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value2 = 1769.6
$ws.Range("I107").Value2 = 1762.25
$ws.Range("J107").Value2 = 1799
$ws.Range("K107").Value2 = 5286.75
$ws.Range("L107").Value2 = 5397
$ws.Range("M107").Value2 = -3366.75
$ws.Range("N107").Value2 = -9237

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value2 = 92498.75
$ws.Range("I69").Value2 = 0
$ws.Range("J69").Value2 = 92498.75
$ws.Range("K69").Value2 = 0
$ws.Range("L69").Value2 = 92498.75
$ws.Range("N69").Value2 = -93996.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H72").Value2 = 92498.75
$ws.Range("I72").Value2 = 0
$ws.Range("J72").Value2 = 92498.75
$ws.Range("K72").Value2 = 0
$ws.Range("L72").Value2 = 277496.25
$ws.Range("N72").Value2 = -284984.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value2 = 17000
$ws.Range("I98").Value2 = 0
$ws.Range("J98").Value2 = 17000
$ws.Range("K98").Value2 = 0
$ws.Range("L98").Value2 = 17000
$ws.Range("N98").Value2 = -22990

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 1261.5294
$ws.Range("I102").Value2 = 889
$ws.Range("J102").Value2 = 3000
$ws.Range("K102").Value2 = 889
$ws.Range("L102").Value2 = 3000
$ws.Range("M102").Value2 = 733
$ws.Range("N102").Value2 = -6244

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value2 = 1958.8889
$ws.Range("I113").Value2 = 1958.75
$ws.Range("J113").Value2 = 1960
$ws.Range("K113").Value2 = 1958.75
$ws.Range("L113").Value2 = 1960
$ws.Range("M113").Value2 = 211.25
$ws.Range("N113").Value2 = -6300

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value2 = 3005.2
$ws.Range("I12").Value2 = 30
$ws.Range("J12").Value2 = 3749
$ws.Range("K12").Value2 = 30
$ws.Range("L12").Value2 = 3749
$ws.Range("M12").Value2 = 140
$ws.Range("N12").Value2 = -4089

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value2 = 4685.6665
$ws.Range("I13").Value2 = 60
$ws.Range("J13").Value2 = 6998.5
$ws.Range("K13").Value2 = 60
$ws.Range("L13").Value2 = 6998.5
$ws.Range("M13").Value2 = 80
$ws.Range("N13").Value2 = -7278.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value2 = 500
$ws.Range("I17").Value2 = 0
$ws.Range("J17").Value2 = 500
$ws.Range("K17").Value2 = 0
$ws.Range("L17").Value2 = 500
$ws.Range("N17").Value2 = -840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value2 = 5502
$ws.Range("I19").Value2 = 0
$ws.Range("J19").Value2 = 5502
$ws.Range("K19").Value2 = 0
$ws.Range("L19").Value2 = 5502
$ws.Range("N19").Value2 = -5842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 277.5
$ws.Range("I22").Value2 = 554
$ws.Range("J22").Value2 = 1
$ws.Range("K22").Value2 = 554
$ws.Range("L22").Value2 = 1
$ws.Range("M22").Value2 = -259
$ws.Range("N22").Value2 = -591

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value2 = 277.5
$ws.Range("I27").Value2 = 554
$ws.Range("J27").Value2 = 1
$ws.Range("K27").Value2 = 554
$ws.Range("L27").Value2 = 1
$ws.Range("M27").Value2 = -447
$ws.Range("N27").Value2 = -215

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value2 = 30103
$ws.Range("I58").Value2 = 0
$ws.Range("J58").Value2 = 30103
$ws.Range("K58").Value2 = 0
$ws.Range("L58").Value2 = 30103
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value2 = -30623

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H97").Value2 = 6000
$ws.Range("I97").Value2 = 0
$ws.Range("J97").Value2 = 6000
$ws.Range("K97").Value2 = 0
$ws.Range("L97").Value2 = 6000
$ws.Range("N97").Value2 = -7982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value2 = 3536.7058
$ws.Range("I122").Value2 = 3317.2307
$ws.Range("J122").Value2 = 4250
$ws.Range("K122").Value2 = 9951.6921
$ws.Range("L122").Value2 = 12750
$ws.Range("M122").Value2 = -7501.6921
$ws.Range("N122").Value2 = -17650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value2 = 20181.5
$ws.Range("I41").Value2 = 19565
$ws.Range("J41").Value2 = 20798
$ws.Range("K41").Value2 = 19565
$ws.Range("L41").Value2 = 20798
$ws.Range("M41").Value2 = -19175
$ws.Range("N41").Value2 = -21578

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value2 = 21694.6
$ws.Range("I74").Value2 = 20401
$ws.Range("J74").Value2 = 22018
$ws.Range("K74").Value2 = 20401
$ws.Range("L74").Value2 = 22018
$ws.Range("M74").Value2 = -19465
$ws.Range("N74").Value2 = -23890

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value2 = 21694.6
$ws.Range("I77").Value2 = 20401
$ws.Range("J77").Value2 = 22018
$ws.Range("K77").Value2 = 61203
$ws.Range("L77").Value2 = 66054
$ws.Range("M77").Value2 = -56523
$ws.Range("N77").Value2 = -75414

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 445.875
$ws.Range("I81").Value2 = 366.7143
$ws.Range("J81").Value2 = 1000
$ws.Range("K81").Value2 = 733.4286
$ws.Range("L81").Value2 = 2000
$ws.Range("M81").Value2 = 327.5714
$ws.Range("N81").Value2 = -4122

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value2 = 445.875
$ws.Range("I84").Value2 = 366.7143
$ws.Range("J84").Value2 = 1000
$ws.Range("K84").Value2 = 3667.143
$ws.Range("L84").Value2 = 10000
$ws.Range("M84").Value2 = 1636.857
$ws.Range("N84").Value2 = -20608

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value2 = 44666.332
$ws.Range("I92").Value2 = 0
$ws.Range("J92").Value2 = 44666.332
$ws.Range("K92").Value2 = 0
$ws.Range("L92").Value2 = 44666.332
$ws.Range("N92").Value2 = -49658.332
